# The document's customXml parts were re-numbered by Word: the SharePoint
# "FormTemplates" content-type part (ds:itemID {1EBB0D41-FC10-49BF-9DA3-5C870D105EC2})
# used to be stored first (customXml/item1.xml + itemProps1.xml) and the
# "contentTypeSchema" part (ds:itemID {9B471222-C2CE-47B9-A9F3-2DFEC7520ECE})
# used to be stored second (customXml/item2.xml + itemProps2.xml). After the
# edit their storage order is swapped: contentTypeSchema becomes item1 /
# itemProps1 and FormTemplates becomes item2 / itemProps2. Their contents are
# otherwise unchanged.
#
# Word re-numbers customXml/item*.xml + itemProps*.xml purely by the order the
# parts are enumerated in Document.CustomXMLParts at save time, so we
# reproduce the new order by removing both parts and re-adding them (via the
# CustomXMLParts collection) in the desired final order: the contentType
# schema part first, then the SharePoint form-templates part.

$d = $word.ActiveDocument
$parts = $d.CustomXMLParts

$formNamespace = "http://schemas.microsoft.com/sharepoint/v3/contenttype/forms"
$schemaNamespace = "http://schemas.microsoft.com/office/2006/metadata/contentType"

$formXml = $null
$schemaXml = $null
$formPart = $null
$schemaPart = $null

for ($i = 1; $i -le $parts.Count; $i++) {
    $p = $parts.Item($i)
    if ($p.NamespaceURI -eq $formNamespace) {
        $formXml = $p.XML
        $formPart = $p
    } elseif ($p.NamespaceURI -eq $schemaNamespace) {
        $schemaXml = $p.XML
        $schemaPart = $p
    }
}

if ($formXml -ne $null -and $schemaXml -ne $null) {
    # Remove the existing copies before re-inserting them in the swapped order.
    $schemaPart.Delete()
    $formPart.Delete()

    # Re-add in the new desired physical order: contentTypeSchema (-> item1 /
    # itemProps1), then FormTemplates (-> item2 / itemProps2).
    $d.CustomXMLParts.Add($schemaXml)
    $d.CustomXMLParts.Add($formXml)
}
